# Preliminary check-in: add a "properties" sheet (table_id config props) so that
# a properties.csv can be generated for this ODK Survey form, minimized to the
# "Table"/"default" partition/aspect entries: colOrder, defaultViewType and
# listViewFileName.

$wb = $excel.ActiveWorkbook

# Append a new worksheet named "properties" after the last existing sheet
# (survey, choices, queries, settings -> ..., properties).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$propsSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$propsSheet.Name = "properties"

# Header row.
$propsSheet.Cells.Item(1, 1).Value = "partition"
$propsSheet.Cells.Item(1, 2).Value = "aspect"
$propsSheet.Cells.Item(1, 3).Value = "key"
$propsSheet.Cells.Item(1, 4).Value = "type"
$propsSheet.Cells.Item(1, 5).Value = "value"

# colOrder (array) - value (the JSON list of columns) is filled in afterwards.
$propsSheet.Cells.Item(2, 1).Value = "Table"
$propsSheet.Cells.Item(2, 2).Value = "default"
$propsSheet.Cells.Item(2, 3).Value = "colOrder"
$propsSheet.Cells.Item(2, 4).Value = "array"

# defaultViewType (string) = LIST
$propsSheet.Cells.Item(3, 1).Value = "Table"
$propsSheet.Cells.Item(3, 2).Value = "default"
$propsSheet.Cells.Item(3, 3).Value = "defaultViewType"
$propsSheet.Cells.Item(3, 4).Value = "string"
$propsSheet.Cells.Item(3, 5).Value = "LIST"

# listViewFileName (configpath)
$propsSheet.Cells.Item(4, 1).Value = "Table"
$propsSheet.Cells.Item(4, 2).Value = "default"
$propsSheet.Cells.Item(4, 3).Value = "listViewFileName"
$propsSheet.Cells.Item(4, 4).Value = "configpath"
$propsSheet.Cells.Item(4, 5).Value = "config/tables/visit/html/visit_list.html"

# colOrder value - the ordered list of columns shown in the table's list view.
$propsSheet.Cells.Item(2, 5).Value = '["plot_id","date","plant_height","plant_health","plant_picture_uriFragment","pests","soil","observations"]'

# Leave the cursor where the author last left it on the new sheet.
$propsSheet.Range("B21").Select() | Out-Null
